$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Z7").Value = "Supported"
$ws.Range("Z8").Value = "Supported"
$ws.Range("Z9").Value = "Supported"
$ws.Range("Z10").Value = "Supported"
$ws.Range("Z11").Value = "Supported"
$ws.Range("Z12").Value = "Supported"
$ws.Range("Z13").Value = "Supported"
$ws.Range("Z14").Value = "Supported"
$ws.Range("Z15").Value = "Supported"
$ws.Range("Z16").Value = "Supported"
$ws.Range("X20").Value = 11000
$ws.Range("X21").Value = 22000
$ws.Range("X27").Value = "Not Supported"
$ws.Range("Y27").Value = "Not Supported"
$ws.Range("X29").Value = 1500
$ws.Range("Y29").Value = 12000
$ws.Range("X30").Value = 3000
$ws.Range("Y30").Value = 25000
$ws.Range("X31").Value = 6000
$ws.Range("Y31").Value = 50000
$ws.Range("X32").Value = 11000
$ws.Range("X33").Value = "Not Supported"
$ws.Range("Y33").Value = "Not Supported"
$ws.Range("X35").Value = 22000
$ws.Range("Y35").Value = "Not Supported"
$ws.Range("X38").Value = 2325
$ws.Range("Y38").Value = "Request"
$ws.Range("X40").Value = 4650
$ws.Range("Y40").Value = "Request"
$ws.Range("X41").Value = "Not Supported"
$ws.Range("X42").Value = 9300
$ws.Range("X44").Value = 18600
$ws.Range("X56").Value = 2325
$ws.Range("Y56").Value = "Request"
$ws.Range("X57").Value = 3530
$ws.Range("Y57").Value = "Request"
$ws.Range("X58").Value = 4650
$ws.Range("Y58").Value = 48750
$ws.Range("X59").Value = 6680
$ws.Range("Y59").Value = "Request"
$ws.Range("X60").Value = 9300
$ws.Range("Y60").Value = 91050
$ws.Range("X61").Value = 12300
$ws.Range("Y61").Value = "Request"
$ws.Range("X62").Value = 18600
$ws.Range("Y62").Value = "Request"
$ws.Range("X63").Value = 24180
$ws.Range("Y63").Value = "Request"
$ws.Range("X64").Value = 30430
$ws.Range("Y64").Value = "Request"
$ws.Range("X83").Value = "Not Supported"
$ws.Range("X85").Value = "Not Supported"
$ws.Range("X87").Value = 3580
$ws.Range("Y87").Value = 34415
$ws.Range("X88").Value = 6900
$ws.Range("Y88").Value = 78620
$ws.Range("X89").Value = 11870
$ws.Range("Y89").Value = 137520
$ws.Range("X90").Value = 22680
$ws.Range("Y90").Value = 247880
$ws.Range("X91").Value = 41670
$ws.Range("Y91").Value = "Request"
$ws.Range("W102").Value = "Unknown"
$ws.Range("X102").Value = 134630
$ws.Range("Y102").Value = "Request"
$ws.Range("Z102").Value = "Not Supported"
$ws.Range("W103").Value = "Unknown"
$ws.Range("X103").Value = "Not Supported"
$ws.Range("Y103").Value = "Not Supported"
$ws.Range("Z103").Value = "Not Supported"
$ws.Range("W104").Value = "Unknown"
$ws.Range("X104").Value = 68930
$ws.Range("Y104").Value = "Request"
$ws.Range("Z104").Value = "Not Supported"
$ws.Range("W105").Value = "Unknown"
$ws.Range("X105").Value = "Not Supported"
$ws.Range("Y105").Value = "Not Supported"
$ws.Range("Z105").Value = "Not Supported"
$ws.Range("W106").Value = "Unknown"
$ws.Range("X106").Value = "Not Supported"
$ws.Range("Y106").Value = "Not Supported"
$ws.Range("Z106").Value = "Not Supported"
$ws.Range("W107").Value = "Unknown"
$ws.Range("X107").Value = "Not Supported"
$ws.Range("Y107").Value = "Not Supported"
$ws.Range("Z107").Value = "Not Supported"
$ws.Range("W108").Value = "Unknown"
$ws.Range("X108").Value = "Not Supported"
$ws.Range("Y108").Value = "Not Supported"
$ws.Range("Z108").Value = "Not Supported"
$ws.Range("W109").Value = "Unknown"
$ws.Range("X109").Value = "Not Supported"
$ws.Range("Y109").Value = "Not Supported"
$ws.Range("Z109").Value = "Not Supported"
$ws.Range("W110").Value = "Unknown"
$ws.Range("X110").Value = "Not Supported"
$ws.Range("Y110").Value = "Not Supported"
$ws.Range("Z110").Value = "Not Supported"
$ws.Range("W111").Value = "Unknown"
$ws.Range("X111").Value = "Not Supported"
$ws.Range("Y111").Value = "Not Supported"
$ws.Range("Z111").Value = "Not Supported"

# Rows 102-111: formulas in W/X/Y were replaced with literal values above;
# restore the plain 'General' number format (style index 1) instead of the
# inherited comma-style format (style index 2) left over from the formulas.
$ws.Range("W102").NumberFormat = "general"
$ws.Range("X102").NumberFormat = "general"
$ws.Range("Y102").NumberFormat = "general"
$ws.Range("W103").NumberFormat = "general"
$ws.Range("X103").NumberFormat = "general"
$ws.Range("Y103").NumberFormat = "general"
$ws.Range("W104").NumberFormat = "general"
$ws.Range("X104").NumberFormat = "general"
$ws.Range("Y104").NumberFormat = "general"
$ws.Range("W105").NumberFormat = "general"
$ws.Range("X105").NumberFormat = "general"
$ws.Range("Y105").NumberFormat = "general"
$ws.Range("W106").NumberFormat = "general"
$ws.Range("X106").NumberFormat = "general"
$ws.Range("Y106").NumberFormat = "general"
$ws.Range("W107").NumberFormat = "general"
$ws.Range("X107").NumberFormat = "general"
$ws.Range("Y107").NumberFormat = "general"
$ws.Range("W108").NumberFormat = "general"
$ws.Range("X108").NumberFormat = "general"
$ws.Range("Y108").NumberFormat = "general"
$ws.Range("W109").NumberFormat = "general"
$ws.Range("X109").NumberFormat = "general"
$ws.Range("Y109").NumberFormat = "general"
$ws.Range("W110").NumberFormat = "general"
$ws.Range("X110").NumberFormat = "general"
$ws.Range("Y110").NumberFormat = "general"
$ws.Range("W111").NumberFormat = "general"
$ws.Range("X111").NumberFormat = "general"
$ws.Range("Y111").NumberFormat = "general"
